{"js": "// Applies the \"Se agreg\u00f3 informaci\u00f3n al informe.\" edit:\n// rewrites several paragraphs in the \"Dise\u00f1o del Framework\" / \"Falencias\n// del Dise\u00f1o\" sections, and adds one new paragraph.\n\nasync function replaceExact(oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText.substring(0, 60));\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"Se comenz\u00f3 por crear...\" paragraph -> new refactor / Composite explanation.\nawait replaceExact(\n  \"Se comenz\u00f3 por crear una relaci\u00f3n contenedor \u2013 contenido entre TestCollection y UnitTest, teniendo TestCollection m\u00e9todos para agregar, correr todos y mostrar el resultado de los tests contenidos. Similarmente se cre\u00f3 otra relaci\u00f3n contenedor \u2013 contenido entre TestReport y TestResult para almacenar resultados de un conjunto de tests, y a TestReport se lo compuso con TestCollection para que est\u00e9 \u00faltimo se abstraiga de la forma en que se almacenan y muestran los resultados.\",\n  \"Partiendo de la base del dise\u00f1o armado para la primer entrega se hizo un refactor de la relaci\u00f3n contenedor \u2013 contenido entre TestCollection y UnitTest, utilizando el patr\u00f3n Composite. Al utilizar este patr\u00f3n los tests se ejecutan en forma recursiva, el usuario s\u00f3lo ejecuta el run del objeto contenedor principal de tipo TestCollection, lo que desencadena la ejecuci\u00f3n de todos los run de los TestCollection y UnitTest contenidos en \u00e9l. La necesidad de untilizar este patr\u00f3n se nos present\u00f3 con el nuevo requerimiento de que se permita almacenar TestCollections dentro de TestCollections. El mismo patr\u00f3n se utiliz\u00f3 para generar la estructura d\u00f3nde se almacenan los resultados de los test.\"\n);\n\n// 2) \"Luego se implement\u00f3...\" -> \"Se mantuvo...\"\nawait replaceExact(\n  \"Luego se implement\u00f3 la clase utilitaria Validation,\",\n  \"Se mantuvo la clase utilitaria Validation,\"\n);\n\n// 3) \"El problema principal a resolver...\" paragraph -> moved/new text about\n//    heredar de UnitTest (duplicated further down, per the diff).\nawait replaceExact(\n  \"El problema principal a resolver fue c\u00f3mo correr un m\u00e9todo creado por el cliente (su m\u00e9todo que hace el testing de algo). Si el lenguaje permitiese pasar m\u00e9todos como argumentos, el problema se resolver\u00eda pas\u00e1ndole dicho m\u00e9todo a UnitTest para que lo ejecute y verifique si alg\u00fan validation fall\u00f3. Como esto no es posible en Java, y al no poder usar reflection, se pens\u00f3 en utilizar la herencia. \",\n  \"Al igual que para las entregas anteriores, para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\"\n);\n\n// 4) \"Para hacer un test...\" paragraph -> new text about TestResultCollection / TestReport.\nawait replaceExact(\n  \"Para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\",\n  \"En el momento en que se ejecutan los tests se imprime por pantalla el resultado de las mismas, y se va generando una estructura de tipo TestResultCollection donde se almacenan objetos de tipo TestResult. Para generar el reporte, el usuario debe pasar por par\u00e1metro esa estructura al constructor de la clase TestReport e invocar al m\u00e9todo generarReporteEnArchivo() de esta \u00faltima.\"\n);\n\n// 5) Insert a brand-new paragraph right after the paragraph from step 4,\n//    re-using (duplicating) the original \"Para hacer un test...\" wording.\n{\n  const results = context.document.body.search(\n    \"En el momento en que se ejecutan los tests se imprime por pantalla el resultado de las mismas, y se va generando una estructura de tipo TestResultCollection donde se almacenan objetos de tipo TestResult. Para generar el reporte, el usuario debe pasar por par\u00e1metro esa estructura al constructor de la clase TestReport e invocar al m\u00e9todo generarReporteEnArchivo() de esta \u00faltima.\",\n    { matchCase: true }\n  );\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Anchor paragraph (step 5) not found\");\n  }\n  const anchorPara = results.items[0].paragraphs.getFirst();\n  anchorPara.insertParagraph(\n    \"Al igual que para las entregas anteriores, para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// 6) \"El principal problema del dise\u00f1o...\" -> prefixed with \"Al igual que para las entregas anteriores, \".\nawait replaceExact(\n  \"El principal problema del dise\u00f1o es que el cliente debe crear tantas clases que hereden de UnitTest como tests quiera hacer.\",\n  \"Al igual que para las entregas anteriores, el principal problema del dise\u00f1o es que el cliente debe crear tantas clases que hereden de UnitTest como tests quiera hacer.\"\n);\n\n// 7) \"Por otro lado, y pensando...\" paragraph -> Composite pattern drawback discussion.\nawait replaceExact(\n  \"Por otro lado, y pensando en los objetivos de la materia, parece ser que se podr\u00eda implementar alg\u00fan patr\u00f3n de creaci\u00f3n de instancias (factory, builder). No se encontr\u00f3 la forma de hacerlo para la creaci\u00f3n de los tests del usuario, ya que a priori no se conocen los nombres que le pondr\u00e1 el usuario a las clases de sus tests (no son conocidos los nombres ni la cantidad).\",\n  \"Con respecto al uso del patr\u00f3n Composite, detectamos que la desventaja de usarlo es que tanto en la clase UnitTest como en UnitTestResult es que nos vimos obligados a implementar m\u00e9todos que no se van a utilizar (como es el caso del m\u00e9todo \\u201cadd\\u201d). Sin embargo, consideramos que era mayor el beneficio obtenido con respecto a la ejecuci\u00f3n de las pruebas, en forma polim\u00f3rfica.\"\n);\n", "ps1": "# Applies the \"Se agreg\u00f3 informaci\u00f3n al informe.\" edit:\n# rewrites several paragraphs in the \"Dise\u00f1o del Framework\" / \"Falencias\n# del Dise\u00f1o\" sections, and adds one new paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Exact {\n    param($Doc, [string]$OldText, [string]$NewText)\n\n    $find = $Doc.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)\n    if (-not $found) {\n        throw \"Replace-Exact: text not found: $OldText\"\n    }\n}\n\n# 1) \"Se comenz\u00f3 por crear...\" paragraph -> new refactor / Composite explanation.\nReplace-Exact $d `\n    \"Se comenz\u00f3 por crear una relaci\u00f3n contenedor \u2013 contenido entre TestCollection y UnitTest, teniendo TestCollection m\u00e9todos para agregar, correr todos y mostrar el resultado de los tests contenidos. Similarmente se cre\u00f3 otra relaci\u00f3n contenedor \u2013 contenido entre TestReport y TestResult para almacenar resultados de un conjunto de tests, y a TestReport se lo compuso con TestCollection para que est\u00e9 \u00faltimo se abstraiga de la forma en que se almacenan y muestran los resultados.\" `\n    \"Partiendo de la base del dise\u00f1o armado para la primer entrega se hizo un refactor de la relaci\u00f3n contenedor \u2013 contenido entre TestCollection y UnitTest, utilizando el patr\u00f3n Composite. Al utilizar este patr\u00f3n los tests se ejecutan en forma recursiva, el usuario s\u00f3lo ejecuta el run del objeto contenedor principal de tipo TestCollection, lo que desencadena la ejecuci\u00f3n de todos los run de los TestCollection y UnitTest contenidos en \u00e9l. La necesidad de untilizar este patr\u00f3n se nos present\u00f3 con el nuevo requerimiento de que se permita almacenar TestCollections dentro de TestCollections. El mismo patr\u00f3n se utiliz\u00f3 para generar la estructura d\u00f3nde se almacenan los resultados de los test.\"\n\n# 2) \"Luego se implement\u00f3...\" -> \"Se mantuvo...\"\nReplace-Exact $d `\n    \"Luego se implement\u00f3 la clase utilitaria Validation,\" `\n    \"Se mantuvo la clase utilitaria Validation,\"\n\n# 3) \"El problema principal a resolver...\" paragraph -> moved/new text about\n#    heredar de UnitTest (duplicated further down, per the diff).\nReplace-Exact $d `\n    \"El problema principal a resolver fue c\u00f3mo correr un m\u00e9todo creado por el cliente (su m\u00e9todo que hace el testing de algo). Si el lenguaje permitiese pasar m\u00e9todos como argumentos, el problema se resolver\u00eda pas\u00e1ndole dicho m\u00e9todo a UnitTest para que lo ejecute y verifique si alg\u00fan validation fall\u00f3. Como esto no es posible en Java, y al no poder usar reflection, se pens\u00f3 en utilizar la herencia. \" `\n    \"Al igual que para las entregas anteriores, para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\"\n\n# 4) \"Para hacer un test...\" paragraph -> new text about TestResultCollection / TestReport.\nReplace-Exact $d `\n    \"Para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\" `\n    \"En el momento en que se ejecutan los tests se imprime por pantalla el resultado de las mismas, y se va generando una estructura de tipo TestResultCollection donde se almacenan objetos de tipo TestResult. Para generar el reporte, el usuario debe pasar por par\u00e1metro esa estructura al constructor de la clase TestReport e invocar al m\u00e9todo generarReporteEnArchivo() de esta \u00faltima.\"\n\n# 5) Insert a brand-new paragraph right after the paragraph from step 4,\n#    re-using (duplicating) the original \"Para hacer un test...\" wording.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$found = $find.Execute(\"En el momento en que se ejecutan los tests se imprime por pantalla el resultado de las mismas, y se va generando una estructura de tipo TestResultCollection donde se almacenan objetos de tipo TestResult. Para generar el reporte, el usuario debe pasar por par\u00e1metro esa estructura al constructor de la clase TestReport e invocar al m\u00e9todo generarReporteEnArchivo() de esta \u00faltima.\")\nif (-not $found) {\n    throw \"Anchor paragraph (step 5) not found\"\n}\n$anchorPara = $find.Parent.Paragraphs(1)\n$anchorPara.Range.InsertParagraphAfter()\n$newPara = $anchorPara.Next()\n$newPara.Range.Text = \"Al igual que para las entregas anteriores, para hacer un test, el usuario debe heredar de UnitTest y redefinir el m\u00e9todo test(), en donde escribe su c\u00f3digo de prueba. El framework se encarga de ejecutarlo e informar el resultado.\"\n\n# 6) \"El principal problema del dise\u00f1o...\" -> prefixed with \"Al igual que para las entregas anteriores, \".\nReplace-Exact $d `\n    \"El principal problema del dise\u00f1o es que el cliente debe crear tantas clases que hereden de UnitTest como tests quiera hacer.\" `\n    \"Al igual que para las entregas anteriores, el principal problema del dise\u00f1o es que el cliente debe crear tantas clases que hereden de UnitTest como tests quiera hacer.\"\n\n# 7) \"Por otro lado, y pensando...\" paragraph -> Composite pattern drawback discussion.\nReplace-Exact $d `\n    \"Por otro lado, y pensando en los objetivos de la materia, parece ser que se podr\u00eda implementar alg\u00fan patr\u00f3n de creaci\u00f3n de instancias (factory, builder). No se encontr\u00f3 la forma de hacerlo para la creaci\u00f3n de los tests del usuario, ya que a priori no se conocen los nombres que le pondr\u00e1 el usuario a las clases de sus tests (no son conocidos los nombres ni la cantidad).\" `\n    \"Con respecto al uso del patr\u00f3n Composite, detectamos que la desventaja de usarlo es que tanto en la clase UnitTest como en UnitTestResult es que nos vimos obligados a implementar m\u00e9todos que no se van a utilizar (como es el caso del m\u00e9todo \u201cadd\u201d). Sin embargo, consideramos que era mayor el beneficio obtenido con respecto a la ejecuci\u00f3n de las pruebas, en forma polim\u00f3rfica.\"\n"}
